# Update "last_edited_time" (column D) for rows 7-12 from 2024-07-19T12:51:00.000Z
# to 2024-07-20T13:34:00.000Z, and refresh the associated numeric rollup values
# in row 7 to match the new snapshot (accounts for Mac/Win calc differences).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("last_edited_time") for rows 7 through 12
foreach ($r in 7..12) {
    $ws.Cells.Item($r, 4).Value = "2024-07-20T13:34:00.000Z"
}

# Row 7 numeric updates
$ws.Range("W7").Value = 224378000
$ws.Range("AA7").Value = 169440000
$ws.Range("AE7").Value = 393818000
$ws.Range("AH7").Value = 329818000
$ws.Range("AK7").Value = 56
$ws.Range("AN7").Value = 64000000
$ws.Range("AQ7").Value = 367118000
